{"js": "// The document has two table-caption bookmarks (\"tab:OverviewTable\" and\n// \"tab:StateLevelTable\") whose internal w:id changed (re-minted) while\n// everything else about them (name, location, surrounding text) stayed the\n// same -- i.e. the bookmarks were re-created in place. Reproduce that by\n// capturing each bookmark's range, deleting the bookmark, and re-inserting\n// a bookmark with the same name over the same range.\n\nconst body = context.document.body;\n\nconst bookmarkNames = [\"tab:OverviewTable\", \"tab:StateLevelTable\"];\nconst ranges = bookmarkNames.map((name) => body.getBookmarkRange(name));\nranges.forEach((r) => r.load(\"text\"));\nawait context.sync();\n\nfor (let i = 0; i < bookmarkNames.length; i++) {\n  const name = bookmarkNames[i];\n  const range = ranges[i];\n  context.document.deleteBookmark(name);\n  range.insertBookmark(name);\n}\n\nawait context.sync();\n", "ps1": "# The document has two table-caption bookmarks (\"tab:OverviewTable\" and\n# \"tab:StateLevelTable\") whose internal w:id changed (re-minted) while\n# everything else about them (name, location, surrounding text) stayed the\n# same -- i.e. the bookmarks were re-created in place. Reproduce that by\n# re-adding each bookmark over its own existing range: Bookmarks.Add with\n# a name that already exists redefines that bookmark.\n\n$d = $word.ActiveDocument\n\n$names = @(\"tab:OverviewTable\", \"tab:StateLevelTable\")\n\nforeach ($name in $names) {\n    $bm = $d.Bookmarks($name)\n    $r = $bm.Range\n    $d.Bookmarks.Add($name, $r)\n}\n"}
